$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Acknowledgements paragraph: split the run so the proper-noun surnames
#    "Inchley" and "Mabelis" get wrapped in spell-check proofErr markers,
#    exactly as Word's background spell-checker would do after a re-save.
# ---------------------------------------------------------------------------
$ackPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*HBSC Scotland research team*") {
        $ackPara = $d.Paragraphs($i)
        break
    }
}
if ($ackPara -ne $null) {
    $xmlFrag = @'
<w:p w14:paraId="3B3681EC" w14:textId="77777777" w:rsidR="004F5BF2" w:rsidRDefault="00325D4C"><w:r><w:t xml:space="preserve">This survey was conducted as part of the Health Behaviour in School-aged Children (HBSC) Scotland study, led by the MRC/CSO Social and Public Health Sciences Unit, University of Glasgow in collaboration with the School of Medicine, University of St Andrews. The HBSC study is funded by Public Health Scotland. The HBSC Scotland research team includes Dr Jo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Inchley</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, Dorothy Currie, Dr Judith Brown, Judith </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mabelis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and Dr Malachi Willis. We would like to thank all the pupils who took part in the survey and the teachers who supported the pilot.</w:t></w:r></w:p>
'@
    $ackPara.Range.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 2) "Generally, within the same figure..." bullet: wrap "and also" in a
#    grammar-check proofErr marker.
# ---------------------------------------------------------------------------
$genderPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Generally, within the same figure*") {
        $genderPara = $d.Paragraphs($i)
        break
    }
}
if ($genderPara -ne $null) {
    $xmlFrag = @'
<w:p w14:paraId="3B368200" w14:textId="77777777" w:rsidR="004F5BF2" w:rsidRDefault="00325D4C" w:rsidP="00366FFF"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Generally, within the same figure, we provide results by gender (in blue) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>and also</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> by year group (in green)</w:t></w:r></w:p>
'@
    $genderPara.Range.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 3) Table cell "15 year-olds": wrap the run in a grammar-check proofErr
#    marker pair.
# ---------------------------------------------------------------------------
for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $tbl = $d.Tables($t)
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $cell = $null
            try { $cell = $tbl.Cell($r, $c) } catch { $cell = $null }
            if ($cell -ne $null -and $cell.Range.Text -like "*15 year-olds*") {
                $xmlFrag = @'
<w:p w14:paraId="3B36820A" w14:textId="77777777" w:rsidR="004F5BF2" w:rsidRDefault="00325D4C"><w:pPr><w:keepNext/><w:spacing w:before="100" w:after="100" w:line="240" w:lineRule="auto"/><w:ind w:left="100" w:right="100"/><w:jc w:val="right"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/></w:rPr><w:t>15 year-olds</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
                $cell.Range.Paragraphs(1).Range.InsertXML($xmlFrag)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Style changes to support the new chart type (page breaks before the
#    higher-level headings so each new graph starts on its own page, and a
#    larger Heading 3 run size to match).
# ---------------------------------------------------------------------------
$h1 = $d.Styles("Heading1")
$h1.ParagraphFormat.PageBreakBefore = $true

$h2 = $d.Styles("Heading2")
$h2.ParagraphFormat.PageBreakBefore = $false

$h6 = $d.Styles("Heading6")
$h6.ParagraphFormat.PageBreakBefore = $false

$h3c = $d.Styles("Heading3Char")
$h3c.Font.Size = 16
$h3c.Font.SizeBi = 16
